# Update the "展览" (Exhibition) and "全部类型" (All types) sheets to the
# refreshed listing: the two "合肥·CW国潮动漫游戏嘉年华" rows (old rows 2-3)
# were removed, every remaining event row shifted up by two rows, the
# running index in column A was renumbered 1..15, and the "想去人数"
# (want-to-go count) in column F was refreshed for several events.

$wb = $excel.ActiveWorkbook

# New F (想去人数) values for the rows that changed, keyed by the row
# number *after* the two rows were deleted and everything shifted up.
$fUpdates = @{
    3  = 7061
    4  = 4192
    5  = 66
    6  = 159
    7  = 32
    10 = 60
    11 = 52
    12 = 187
    13 = 611
    14 = 90
    15 = 50
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Remove the two obsolete rows (old rows 2 and 3); remaining rows shift up.
    $ws.Range("A2:A3").EntireRow.Delete()

    # Renumber the running index in column A (1..15) for the 15 remaining
    # data rows now sitting in rows 2..16.
    for ($r = 2; $r -le 16; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 1
    }

    # Refresh the "想去人数" (column F) counts that changed.
    foreach ($r in $fUpdates.Keys) {
        $ws.Cells.Item($r, 6).Value = $fUpdates[$r]
    }
}
